$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value2 = 12.85240442524551
$ws.Range("C2").Value2 = 9.0788740838943
$ws.Range("E2").Value2 = 20.62834473828746
$ws.Range("F2").Value2 = 38.92350550351049
$ws.Range("G2").Value2 = 28.42488304902704
$ws.Range("H2").Value2 = 13.95966697617882
$ws.Range("I2").Value2 = 19.47283218963927
$ws.Range("J2").Value2 = 7.796946191828069
$ws.Range("M2").Value2 = 19.21422926399112
$ws.Range("B3").Value2 = 12.22708506964004
$ws.Range("C3").Value2 = 8.510838904662814
$ws.Range("E3").Value2 = 20.60201019245705
$ws.Range("F3").Value2 = 38.87529702873183
$ws.Range("G3").Value2 = 28.38786841492724
$ws.Range("H3").Value2 = 14.02200283053265
$ws.Range("I3").Value2 = 19.6157394895085
$ws.Range("J3").Value2 = 7.825545148011106
$ws.Range("M3").Value2 = 18.99423397864439
$ws.Range("B4").Value2 = 11.82744436561093
$ws.Range("C4").Value2 = 8.14075129727413
$ws.Range("E4").Value2 = 20.5887531096169
$ws.Range("F4").Value2 = 38.85779920985581
$ws.Range("G4").Value2 = 28.38141716824682
$ws.Range("H4").Value2 = 14.06414934920251
$ws.Range("I4").Value2 = 19.7097307757095
$ws.Range("J4").Value2 = 7.843948740452046
$ws.Range("M4").Value2 = 18.86084871705897
$ws.Range("B5").Value2 = 11.66084029080158
$ws.Range("C5").Value2 = 7.984585093640034
$ws.Range("E5").Value2 = 20.58408681918998
$ws.Range("F5").Value2 = 38.85371243521452
$ws.Range("G5").Value2 = 28.38286222507772
$ws.Range("H5").Value2 = 14.08229369517192
$ws.Range("I5").Value2 = 19.74959718781722
$ws.Range("J5").Value2 = 7.851661253850839
$ws.Range("M5").Value2 = 18.80697319909702
$ws.Range("B6").Value2 = 11.63295592773336
$ws.Range("C6").Value2 = 7.958330829450379
$ws.Range("E6").Value2 = 20.58335654355108
$ws.Range("F6").Value2 = 38.85321761772263
$ws.Range("G6").Value2 = 28.383347552362
$ws.Range("H6").Value2 = 14.08536496480388
$ws.Range("I6").Value2 = 19.75631126783288
$ws.Range("J6").Value2 = 7.852954791201823
$ws.Range("M6").Value2 = 18.79805777858418
$ws.Range("B7").Value2 = 11.82521237113485
$ws.Range("C7").Value2 = 8.138666845022145
$ws.Range("E7").Value2 = 20.58868719342323
$ws.Range("F7").Value2 = 38.85773177196213
$ws.Range("G7").Value2 = 28.38142019134299
$ws.Range("H7").Value2 = 14.06439013139686
$ws.Range("I7").Value2 = 19.71026210298942
$ws.Range("J7").Value2 = 7.844051891110006
$ws.Range("M7").Value2 = 18.86012011798859
$ws.Range("B8").Value2 = 12.64016075020841
$ws.Range("C8").Value2 = 8.887441144782317
$ws.Range("E8").Value2 = 20.61866245563203
$ws.Range("F8").Value2 = 38.9043725689946
$ws.Range("G8").Value2 = 28.40873323523803
$ws.Range("H8").Value2 = 13.98035412597317
$ws.Range("I8").Value2 = 19.5208065045801
$ws.Range("J8").Value2 = 7.80663244577042
$ws.Range("M8").Value2 = 19.13805642541863
$ws.Range("B9").Value2 = 14.10654509536595
$ws.Range("C9").Value2 = 10.18664595776804
$ws.Range("E9").Value2 = 20.70037519937905
$ws.Range("F9").Value2 = 39.09170557795041
$ws.Range("G9").Value2 = 28.59192075960036
$ws.Range("H9").Value2 = 13.8464758632442
$ws.Range("I9").Value2 = 19.19912964764541
$ws.Range("J9").Value2 = 7.739913425982957
$ws.Range("M9").Value2 = 19.69410168647617
$ws.Range("B10").Value2 = 15.09574492080702
$ws.Range("C10").Value2 = 11.03827851807474
$ws.Range("E10").Value2 = 20.77413088818764
$ws.Range("F10").Value2 = 39.28737756873548
$ws.Range("G10").Value2 = 28.80568285287261
$ws.Range("H10").Value2 = 13.76722392980784
$ws.Range("I10").Value2 = 18.99358623160506
$ws.Range("J10").Value2 = 7.69490706912853
$ws.Range("M10").Value2 = 20.10608891344972
$ws.Range("B11").Value2 = 15.5253956051199
$ws.Range("C11").Value2 = 11.40346662431985
$ws.Range("E11").Value2 = 20.81060068340708
$ws.Range("F11").Value2 = 39.38884431368134
$ws.Range("G11").Value2 = 28.91999621999826
$ws.Range("H11").Value2 = 13.73537418767391
$ws.Range("I11").Value2 = 18.90685715963226
$ws.Range("J11").Value2 = 7.675293372409449
$ws.Range("M11").Value2 = 20.29357732789955
$ws.Range("B12").Value2 = 15.68508950810631
$ws.Range("C12").Value2 = 11.53856783876933
$ws.Range("E12").Value2 = 20.82482438817546
$ws.Range("F12").Value2 = 39.4290402106588
$ws.Range("G12").Value2 = 28.96571694591065
$ws.Range("H12").Value2 = 13.72392223661954
$ws.Range("I12").Value2 = 18.87499724001075
$ws.Range("J12").Value2 = 7.667989050372928
$ws.Range("M12").Value2 = 20.36452824246698
$ws.Range("B13").Value2 = 15.65083136467606
$ws.Range("C13").Value2 = 11.50961299955963
$ws.Range("E13").Value2 = 20.82174277650447
$ws.Range("F13").Value2 = 39.42030479015296
$ws.Range("G13").Value2 = 28.95576240225212
$ws.Range("H13").Value2 = 13.7263614684767
$ws.Range("I13").Value2 = 18.88181500967656
$ws.Range("J13").Value2 = 7.669556709219282
$ws.Range("M13").Value2 = 20.349250680826
$ws.Range("B14").Value2 = 15.53859436180596
$ws.Range("C14").Value2 = 11.41464531505286
$ws.Range("E14").Value2 = 20.81176262827645
$ws.Range("F14").Value2 = 39.39211584361642
$ws.Range("G14").Value2 = 28.92370908852627
$ws.Range("H14").Value2 = 13.73441980059195
$ws.Range("I14").Value2 = 18.90421628486937
$ws.Range("J14").Value2 = 7.67468998113069
$ws.Range("M14").Value2 = 20.29941576680701
$ws.Range("B15").Value2 = 15.46945241755949
$ws.Range("C15").Value2 = 11.35606001310106
$ws.Range("E15").Value2 = 20.80570314319802
$ws.Range("F15").Value2 = 39.3750795548903
$ws.Range("G15").Value2 = 28.90439154141877
$ws.Range("H15").Value2 = 13.73943518647874
$ws.Range("I15").Value2 = 18.91806591798887
$ws.Range("J15").Value2 = 7.677850250707635
$ws.Range("M15").Value2 = 20.26888265262266
$ws.Range("B16").Value2 = 15.06724970154601
$ws.Range("C16").Value2 = 11.0139660814803
$ws.Range("E16").Value2 = 20.7718056957729
$ws.Range("F16").Value2 = 39.28099560093874
$ws.Range("G16").Value2 = 28.79855400408848
$ws.Range("H16").Value2 = 13.76939033938682
$ws.Range("I16").Value2 = 18.99939120085999
$ws.Range("J16").Value2 = 7.696206114540523
$ws.Range("M16").Value2 = 20.0938326471415
$ws.Range("B17").Value2 = 14.81523940797071
$ws.Range("C17").Value2 = 10.79841582414693
$ws.Range("E17").Value2 = 20.75175377172976
$ws.Range("F17").Value2 = 39.22645600597456
$ws.Range("G17").Value2 = 28.73798514204011
$ws.Range("H17").Value2 = 13.78884645230401
$ws.Range("I17").Value2 = 19.05102246241679
$ws.Range("J17").Value2 = 7.707686592308568
$ws.Range("M17").Value2 = 19.98642273619178
$ws.Range("B18").Value2 = 14.66837890526599
$ws.Range("C18").Value2 = 10.67234461471183
$ws.Range("E18").Value2 = 20.74049548460477
$ws.Range("F18").Value2 = 39.19625985000631
$ws.Range("G18").Value2 = 28.70475621898596
$ws.Range("H18").Value2 = 13.80043237200385
$ws.Range("I18").Value2 = 19.08135644633105
$ws.Range("J18").Value2 = 7.714370840657262
$ws.Range("M18").Value2 = 19.92465367795192
$ws.Range("B19").Value2 = 14.6183289709126
$ws.Range("C19").Value2 = 10.62929947230616
$ws.Range("E19").Value2 = 20.73673104500044
$ws.Range("F19").Value2 = 39.18623802108512
$ws.Range("G19").Value2 = 28.69378232749953
$ws.Range("H19").Value2 = 13.80442291445911
$ws.Range("I19").Value2 = 19.09173617214161
$ws.Range("J19").Value2 = 7.716647942179298
$ws.Range("M19").Value2 = 19.90374326923515
$ws.Range("B20").Value2 = 14.8422647292111
$ws.Range("C20").Value2 = 10.82157801811487
$ws.Range("E20").Value2 = 20.75385991049662
$ws.Range("F20").Value2 = 39.2321405010335
$ws.Range("G20").Value2 = 28.74426644576694
$ws.Range("H20").Value2 = 13.78673437270922
$ws.Range("I20").Value2 = 19.04546023388002
$ws.Range("J20").Value2 = 7.706456099644796
$ws.Range("M20").Value2 = 19.99785605989178
$ws.Range("B21").Value2 = 15.57164322677956
$ws.Range("C21").Value2 = 11.44262606457181
$ws.Range("E21").Value2 = 20.81468287202446
$ws.Range("F21").Value2 = 39.40034767251495
$ws.Range("G21").Value2 = 28.93305811230999
$ws.Range("H21").Value2 = 13.73203631303944
$ws.Range("I21").Value2 = 18.89760975404137
$ws.Range("J21").Value2 = 7.673178883191907
$ws.Range("M21").Value2 = 20.31405521344277
$ws.Range("B22").Value2 = 16.03078853135026
$ws.Range("C22").Value2 = 11.82994181976383
$ws.Range("E22").Value2 = 20.85684056454621
$ws.Range("F22").Value2 = 39.52060248051392
$ws.Range("G22").Value2 = 29.07060893990076
$ws.Range("H22").Value2 = 13.69983889543351
$ws.Range("I22").Value2 = 18.80671176240982
$ws.Range("J22").Value2 = 7.652146678295823
$ws.Range("M22").Value2 = 20.52041060186567
$ws.Range("B23").Value2 = 15.78736236393437
$ws.Range("C23").Value2 = 11.62492068886489
$ws.Range("E23").Value2 = 20.83412217952193
$ws.Range("F23").Value2 = 39.45548262111907
$ws.Range("G23").Value2 = 28.99590858072136
$ws.Range("H23").Value2 = 13.7166968691744
$ws.Range("I23").Value2 = 18.85469850251292
$ws.Range("J23").Value2 = 7.66330664231358
$ws.Range("M23").Value2 = 20.41032087902132
$ws.Range("B24").Value2 = 14.83005273138046
$ws.Range("C24").Value2 = 10.81111308306035
$ws.Range("E24").Value2 = 20.75290688418941
$ws.Range("F24").Value2 = 39.22956692794011
$ws.Range("G24").Value2 = 28.74142170478618
$ws.Range("H24").Value2 = 13.78768799711941
$ws.Range("I24").Value2 = 19.04797289169173
$ws.Range("J24").Value2 = 7.70701214374586
$ws.Range("M24").Value2 = 19.99268710487356
$ws.Range("B25").Value2 = 13.72485182075879
$ws.Range("C25").Value2 = 9.853272583955086
$ws.Range("E25").Value2 = 20.67583982483934
$ws.Range("F25").Value2 = 39.0307929572477
$ws.Range("G25").Value2 = 28.52844150334026
$ws.Range("H25").Value2 = 13.8793565072043
$ws.Range("I25").Value2 = 19.28077446563222
$ws.Range("J25").Value2 = 7.757254647108032
$ws.Range("M25").Value2 = 19.54283384174385
